$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

# Each populated table row is updated in place by setting the text of its
# existing cells (left-to-right). The cell count for every affected row is
# unchanged before/after, so no row/cell structural edits are required -
# this also preserves each cell's existing formatting (font/size/alignment).

# Word table row 1 (grid row 0)
$row = $tbl.Rows.Item(1)
$row.Cells.Item(1).Range.Text = "97÷8=12, 1"
$row.Cells.Item(2).Range.Text = "24÷7=3, 3"
$row.Cells.Item(3).Range.Text = "30÷3=10, 0"
$row.Cells.Item(4).Range.Text = "89÷6=14, 5"
$row.Cells.Item(5).Range.Text = "77÷4=19, 1"

# Word table row 5 (grid row 4)
$row = $tbl.Rows.Item(5)
$row.Cells.Item(1).Range.Text = "41÷2=20, 1"
$row.Cells.Item(2).Range.Text = "12÷5=2, 2"
$row.Cells.Item(3).Range.Text = "59÷9=6, 5"
$row.Cells.Item(4).Range.Text = "26÷7=3, 5"
$row.Cells.Item(5).Range.Text = "49÷3=16, 1"

# Word table row 9 (grid row 8)
$row = $tbl.Rows.Item(9)
$row.Cells.Item(1).Range.Text = "89÷4=22, 1"
$row.Cells.Item(2).Range.Text = "18÷8=2, 2"
$row.Cells.Item(3).Range.Text = "74÷2=37, 0"
$row.Cells.Item(4).Range.Text = "42÷9=4, 6"
$row.Cells.Item(5).Range.Text = "80÷3=26, 2"

# Word table row 13 (grid row 12)
$row = $tbl.Rows.Item(13)
$row.Cells.Item(1).Range.Text = "65÷2=32, 1"
$row.Cells.Item(2).Range.Text = "38÷8=4, 6"
$row.Cells.Item(3).Range.Text = "25÷7=3, 4"
$row.Cells.Item(4).Range.Text = "81÷9=9, 0"
$row.Cells.Item(5).Range.Text = "62÷5=12, 2"

# Word table row 17 (grid row 16)
$row = $tbl.Rows.Item(17)
$row.Cells.Item(1).Range.Text = "72÷2=36, 0"
$row.Cells.Item(2).Range.Text = "81÷8=10, 1"
$row.Cells.Item(3).Range.Text = "77÷7=11, 0"
$row.Cells.Item(4).Range.Text = "56÷5=11, 1"
$row.Cells.Item(5).Range.Text = "34÷5=6, 4"

Write-Output "updated 5 table rows (25 cells)"
